$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: "Technical Proposal" -> "Technical " (trailing space), no longer wraps to 2 lines ---
$ws.Range("A3").Value = "Technical "
$ws.Rows(3).EntireRow.AutoFit()

# --- Column D gets wider to host the new policy text ---
# (target stored width is 49.88671875; this runtime quantizes ColumnWidth
#  writes to its own pixel grid, so 49 is the input that lands closest to it)
$ws.Columns(4).ColumnWidth = 49

# --- New compliance rows 8-13: Criteria name in column A, "Y/N" in column B ---
# Carry over the existing "Criteria" column formatting (vertical-center +
# wrap text) from row 4 instead of re-applying alignment properties by hand,
# so no redundant/unused cell style gets created in the workbook.
$ws.Range("A4").Copy()
$ws.Range("A8:A13").PasteSpecial(-4122)

$newRows = @(
    @{ Row = 8;  Text = "Code of Conduct" },
    @{ Row = 9;  Text = "Insurance" },
    @{ Row = 10; Text = "Contract terms" },
    @{ Row = 11; Text = "Conflict of Interest" },
    @{ Row = 12; Text = "Governance" },
    @{ Row = 13; Text = "Human Rights" }
)

foreach ($item in $newRows) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = $item.Text
    $ws.Cells.Item($r, 2).Value = "Y/N"
}

# Row 11 ("Conflict of Interest") wraps onto two lines at this column width
$ws.Rows(11).RowHeight = 28.8

# --- Final selection, mirroring where the author clicked when done ---
$ws.Range("D16").Select()
